$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1113552.4
$ws.Range("J17").Value = 1113552.4
$ws.Range("L17").Value = 3340657.2
$ws.Range("N17").Value = -3340993.2
$ws.Range("H62").Value = 3671.4614
$ws.Range("I62").Value = 3394.0833
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 3394.0833
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = -2770.0833
$ws.Range("N62").Value = -8248
$ws.Range("H65").Value = 3671.4614
$ws.Range("I65").Value = 3394.0833
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 16970.4165
$ws.Range("L65").Value = 35000
$ws.Range("M65").Value = -13850.4165
$ws.Range("N65").Value = -41240
$ws.Range("H76").Value = 19242.25
$ws.Range("I76").Value = 11234.5
$ws.Range("K76").Value = 11234.5
$ws.Range("M76").Value = -10919.5
$ws.Range("H79").Value = 19242.25
$ws.Range("I79").Value = 11234.5
$ws.Range("K79").Value = 11234.5
$ws.Range("M79").Value = -10142.5
$ws.Range("H86").Value = 9728.429
$ws.Range("I86").Value = 17359.8
$ws.Range("K86").Value = 17359.8
$ws.Range("M86").Value = -16236.8
$ws.Range("H89").Value = 9728.429
$ws.Range("I89").Value = 17359.8
$ws.Range("K89").Value = 86799
$ws.Range("M89").Value = -81183
$ws.Range("H112").Value = 4236.4165
$ws.Range("I112").Value = 4758.3125
$ws.Range("K112").Value = 14274.9375
$ws.Range("M112").Value = -13166.9375
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 25954
$ws.Range("I33").Value = 25954
$ws.Range("K33").Value = 25954
$ws.Range("M33").Value = -25625
$ws.Range("H112").Value = 30249.25
$ws.Range("J112").Value = 30249.25
$ws.Range("L112").Value = 30249.25
$ws.Range("N112").Value = -33203.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 27308.47
$ws.Range("I20").Value = 11287.1875
$ws.Range("J20").Value = 41549.61
$ws.Range("K20").Value = 11287.1875
$ws.Range("L20").Value = 41549.61
$ws.Range("M20").Value = -11040.1875
$ws.Range("N20").Value = -42043.61
$ws.Range("H22").Value = 2169.8
$ws.Range("I22").Value = 2169.8
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2169.8
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1996.8
$ws.Range("N22").ClearContents()
$ws.Range("H105").Value = 1189.4073
$ws.Range("I105").Value = 973.05
$ws.Range("J105").Value = 1807.5714
$ws.Range("K105").Value = 973.05
$ws.Range("L105").Value = 1807.5714
$ws.Range("M105").Value = 773.95
$ws.Range("N105").Value = -5301.5714
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H134").Value = 97506.71000000001
$ws.Range("I134").Value = 97088.836
$ws.Range("K134").Value = 291266.508
$ws.Range("M134").Value = -288731.508
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 16302.111
$ws.Range("I58").Value = 5778.913
$ws.Range("K58").Value = 5778.913
$ws.Range("M58").Value = -5575.913
$ws.Range("H99").Value = 2915.182
$ws.Range("I99").Value = 3007.5
$ws.Range("J99").Value = 2669
$ws.Range("K99").Value = 3007.5
$ws.Range("L99").Value = 2669
$ws.Range("M99").Value = -1509.5
$ws.Range("N99").Value = -5665
$ws.Range("H126").Value = 2915.182
$ws.Range("I126").Value = 3007.5
$ws.Range("J126").Value = 2669
$ws.Range("K126").Value = 9022.5
$ws.Range("L126").Value = 8007
$ws.Range("M126").Value = -6552.5
$ws.Range("N126").Value = -12947
$ws.Range("H136").Value = 16302.111
$ws.Range("I136").Value = 5778.913
$ws.Range("K136").Value = 17336.739
$ws.Range("M136").Value = -14786.739
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 105.71429
$ws.Range("I12").Value = 38
$ws.Range("J12").Value = 275
$ws.Range("K12").Value = 114
$ws.Range("L12").Value = 825
$ws.Range("M12").Value = 59
$ws.Range("N12").Value = -1171
$ws.Range("H132").Value = 2756079
$ws.Range("I132").Value = 1268.6
$ws.Range("J132").Value = 5051754.5
$ws.Range("K132").Value = 11417.4
$ws.Range("L132").Value = 45465790.5
$ws.Range("M132").Value = -8887.4
$ws.Range("N132").Value = -45470850.5
$ws.Range("H138").Value = 3112.5
$ws.Range("I138").Value = 3112.5
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 9337.5
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -4197.5
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 4010.0527
$ws.Range("I140").Value = 2926.6365
$ws.Range("J140").Value = 5499.75
$ws.Range("K140").Value = 8779.9095
$ws.Range("L140").Value = 16499.25
$ws.Range("M140").Value = -3599.9095
$ws.Range("N140").Value = -26859.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1182
$ws.Range("I97").Value = 1042.9
$ws.Range("K97").Value = 1042.9
$ws.Range("M97").Value = -546.9000000000001
$ws.Range("H134").Value = 76149
$ws.Range("J134").Value = 76149
$ws.Range("L134").Value = 228447
$ws.Range("N134").Value = -233517
$ws.Range("H135").Value = 151123.75
$ws.Range("J135").Value = 151123.75
$ws.Range("L135").Value = 151123.75
$ws.Range("N135").Value = -161263.75
$ws.Range("H139").Value = 93000
$ws.Range("J139").Value = 93000
$ws.Range("L139").Value = 93000
$ws.Range("N139").Value = -103280
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3643.6667
$ws.Range("J46").Value = 3643.6667
$ws.Range("L46").Value = 3643.6667
$ws.Range("N46").Value = -4019.6667
$ws.Range("H110").Value = 225000
$ws.Range("J110").Value = 225000
$ws.Range("L110").Value = 225000
$ws.Range("N110").Value = -233180
$ws.Range("H132").Value = 1705905.1
$ws.Range("I132").Value = 2998.6072
$ws.Range("J132").Value = 5373704
$ws.Range("K132").Value = 8995.821599999999
$ws.Range("L132").Value = 16121112
$ws.Range("M132").Value = -6465.821599999999
$ws.Range("N132").Value = -16126172
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 20941
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H122").Value = 3819.3
$ws.Range("I122").Value = 2939.8
$ws.Range("J122").Value = 6457.8
$ws.Range("K122").Value = 8819.400000000001
$ws.Range("L122").Value = 19373.4
$ws.Range("M122").Value = -6369.400000000001
$ws.Range("N122").Value = -24273.4
$ws.Range("H132").Value = 990365.75
$ws.Range("I132").Value = 2445.8
$ws.Range("J132").Value = 10869565
$ws.Range("K132").Value = 7337.400000000001
$ws.Range("L132").Value = 32608695
$ws.Range("M132").Value = -4807.400000000001
$ws.Range("N132").Value = -32613755
$ws.Range("H137").Value = 89825
$ws.Range("I137").Value = 89650
$ws.Range("J137").Value = 90000
$ws.Range("K137").Value = 89650
$ws.Range("L137").Value = 90000
$ws.Range("M137").Value = -84550
$ws.Range("N137").Value = -100200
$ws.Range("H138").Value = 100327
$ws.Range("J138").Value = 100327
$ws.Range("L138").Value = 100327
$ws.Range("N138").Value = -110607

Write-Host "Applied all updates"